$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment H2:H42 by 1
for ($r = 2; $r -le 42; $r++) {
    $cell = $ws.Cells.Item($r, 8)  # column H = 8
    $current = $cell.Value()
    $cell.Value = $current + 1
}

# Clear formulas/values in columns P and Q (rows 2:42), removing the extra helper columns
$ws.Range("P2:Q42").Clear()

# Update the view: scroll so column F is the left-most visible column,
# and select cell R12 (matches the saved view state in the workbook)
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("R12").Select()
